$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("H1").Value = "Median Household Income"
$ws.Range("I1").Value = "SNAP Recipients"

# New column data: Median Household Income (H) and SNAP Recipients (I)
$ws.Range("H2").Value = 80426
$ws.Range("I2").Value = 36551

$ws.Range("H3").Value = 78386
$ws.Range("I3").Value = 37964

$ws.Range("H4").Value = 67810
$ws.Range("I4").Value = 37893

$ws.Range("H5").Value = 68609
$ws.Range("I5").Value = 36578

$ws.Range("H6").Value = 63755
$ws.Range("I6").Value = 37088

$ws.Range("H7").Value = 62750
$ws.Range("I7").Value = 38675

# Match the bestFit-ish column widths from the target workbook as closely as
# this engine's pixel-rounding allows.
$ws.Columns.Item(8).ColumnWidth = 20.8
$ws.Columns.Item(9).ColumnWidth = 12.76

# Update the view: scroll right so column E is the left-most visible column
# and leave the selection on G18, matching the saved view state.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("G18").Select() | Out-Null
